$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new "Wins", "Losses", "Ties" columns, styled like the
# other header cells (copy the format from AC1 so the same style index is reused).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2 through 40: season record (Wins/Losses/Ties) for each player row.
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 101  # AD
    $ws.Cells.Item($r, 31).Value = 61   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
